$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.303.78"
$ws.Range("E2").Value = "  +3.98%  "
$ws.Range("D3").Value = "1.732.19"
$ws.Range("E3").Value = "  +2.65%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'219.40"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.50%  "
$ws.Range("D6").Value = "'0.522"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.12%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("D8").Value = "'24.13"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +4.26%  "
$ws.Range("E9").Value = "  +2.42%  "
$ws.Range("E10").Value = "  +1.33%  "
$ws.Range("E11").Value = "  +0.35%  "
$ws.Range("D12").Value = "1.977.00"
$ws.Range("E12").Value = "  +2.68%  "
$ws.Range("D13").Value = "1.731.50"
$ws.Range("E13").Value = "  +1.97%  "
$ws.Range("D14").Value = "'4.26"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.77%  "
$ws.Range("D15").Value = "'0.564"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.54%  "
$ws.Range("D16").Value = "'67.81"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.78%  "
$ws.Range("D17").Value = "28.309.17"
$ws.Range("E17").Value = "  +4.02%  "
$ws.Range("D18").Value = "'244.15"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.16%  "
$ws.Range("D19").Value = "0.0₃0755"
$ws.Range("E19").Value = "  +1.16%  "
$ws.Range("D20").Value = "'7.96"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.85%  "
$ws.Range("E21").Value = "  -0.10%  "
$ws.Range("E22").Value = "  +1.71%  "
$ws.Range("D23").Value = "'9.72"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.70%  "
$ws.Range("D24").Value = "'2.10"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.38%  "
$ws.Range("D25").Value = "'149.38"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.66%  "
$ws.Range("D26").Value = "'7.54"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.89%  "
$ws.Range("E27").Value = "  +0.86%  "
$ws.Range("E28").Value = "  +0.54%  "
$ws.Range("E29").Value = "  -0.03%  "
$ws.Range("D30").Value = "'0.0517"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.92%  "
$ws.Range("E31").Value = "  +2.39%  "
$ws.Range("E32").Value = "  +0.45%  "
$ws.Range("E33").Value = "  +1.32%  "
$ws.Range("D34").Value = "1.490.87"
$ws.Range("E34").Value = "  -5.39%  "
$ws.Range("E35").Value = "  -1.71%  "
$ws.Range("D36").Value = "'0.977"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.09%  "
$ws.Range("E37").Value = "  -0.13%  "
$ws.Range("E38").Value = "  +0.88%  "
$ws.Range("E39").Value = "  +0.89%  "
$ws.Range("E40").Value = "  +0.39%  "
$ws.Range("D41").Value = "'70.26"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.69%  "
$ws.Range("E42").Value = "  -0.09%  "
$ws.Range("D43").Value = "'5.67"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.39%  "
$ws.Range("E44").Value = "  +2.17%  "
$ws.Range("E45").Value = "  +2.41%  "
$ws.Range("E46").Value = "  +1.12%  "
$ws.Range("E47").Value = "  +7.40%  "
$ws.Range("E48").Value = "  +5.24%  "
$ws.Range("D49").Value = "'90.82"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.85%  "
$ws.Range("D50").Value = "'8.20"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.43%  "
$ws.Range("D51").Value = "'0.104"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.97%  "
